# Tax, Discount and Price list added to sales module
# Rework PurchaseBillHeader (sheet1) and PurchaseBillItems (sheet2)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: PurchaseBillHeader
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("PurchaseBillHeader")

# Header row styling: give Entry/Invoice/Expected-Delivery-Date columns a
# text number format (keeps their existing bold font + alignment).
$ws1.Range("B1").NumberFormat = "@"
$ws1.Range("C1:D1").NumberFormat = "@"

# Row 2 becomes the single data row, with new date-as-text values.
$ws1.Range("A2").Value = "Automated one"
$ws1.Range("B2").NumberFormat = "@"
$ws1.Range("B2").Value = "23-12-2025"
$ws1.Range("C2").Value = "24-12-2025"
$ws1.Range("D2").NumberFormat = "@"
$ws1.Range("D2").Value = "30-12-2025"
$ws1.Range("E2").Value = "Test S2"
$ws1.Range("F2").Value = "Net 30"

# Remove the old second data row (previously "Automated one") entirely -
# shifts remaining rows up.
$ws1.Rows("3").Delete()

# Drop the long tail of blank placeholder rows (old rows 6-39, now 5-38).
$ws1.Rows("5:38").Delete()

# The two remaining blank rows (3 and 4) should only carry A/B formatting,
# matching the target - clear the stray formatted C cells.
$ws1.Range("C3:C4").Clear()

[void]$ws1.Range("C25").Select()

# ---------------------------------------------------------------------
# Sheet 2: PurchaseBillItems
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("PurchaseBillItems")

# Remove the two "ReTest20" rows - the old row 4 ("Automated one" /
# "BP Apparatus") shifts up to become the new row 2.
$ws2.Rows("2:3").Delete()

# Update the quantity on the remaining data row.
$ws2.Range("C2").Value = 20

[void]$ws2.Range("C2").Select()

[void]$ws1.Select()
